$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-13 (Generation 0-11): Fitness 7573 -> 7598
$ws.Range("C2:C13").Value = 7598

# Rows 14-252 (Generation 12-250): Fitness 7573 -> 7590
$ws.Range("C14:C252").Value = 7590
